$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells stay text so values like "28.696.30" are not
# reinterpreted as numbers/dates by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.696.30"
$ws.Range("E2").Value = "  -1.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.801.81"
$ws.Range("E3").Value = "  -1.49%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.23"
$ws.Range("E5").Value = "  -2.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5953"
$ws.Range("E6").Value = "  -2.79%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("E8").Value = "  -1.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06826"
$ws.Range("E9").Value = "  -4.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.34"
$ws.Range("E10").Value = "  -1.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07544"
$ws.Range("E11").Value = "  -1.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.802.29"
$ws.Range("E12").Value = "  -1.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.760"
$ws.Range("E13").Value = "  -1.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6234"
$ws.Range("E14").Value = "  -1.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.047.18"
$ws.Range("E15").Value = "  -1.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009314"
$ws.Range("E16").Value = "  -7.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "75.39"
$ws.Range("E17").Value = "  -4.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.661.94"
$ws.Range("E18").Value = "  -1.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.469"
$ws.Range("E19").Value = "  -6.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "210.22"
$ws.Range("E21").Value = "  -7.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.44"
$ws.Range("E22").Value = "  -3.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.843"
$ws.Range("E23").Value = "  -2.52%  "

$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.44"
$ws.Range("E25").Value = "  -0.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.841"
$ws.Range("E26").Value = "  -2.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1274"
$ws.Range("E27").Value = "  -3.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.38"
$ws.Range("E28").Value = "  -1.46%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.426"
$ws.Range("E29").Value = "  -3.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06147"
$ws.Range("E30").Value = "  -3.88%  "

$ws.Range("E31").Value = "  -2.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.782"
$ws.Range("E32").Value = "  -1.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.746"
$ws.Range("E33").Value = "  -1.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.718"
$ws.Range("E34").Value = "  -1.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.062"
$ws.Range("E35").Value = "  -6.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6385"
$ws.Range("E36").Value = "  -1.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.489"
$ws.Range("E37").Value = "  -2.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.713"
$ws.Range("E38").Value = "  -1.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.130.95"
$ws.Range("E41").Value = "  -7.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8701"
$ws.Range("E42").Value = "  -5.26%  "

$ws.Range("E43").Value = "  +0.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.65"
$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.965.67"
$ws.Range("E45").Value = "  -0.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.52"
$ws.Range("E46").Value = "  -3.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000112"
$ws.Range("E47").Value = "  -4.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.602"
$ws.Range("E48").Value = "  -1.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05475"
$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.339"
$ws.Range("E50").Value = "  -3.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4490"
$ws.Range("E51").Value = "  -1.77%  "

# Rows 39 and 40 swap coin identity (VeChain <-> FraxShare) along with new values
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.452"
$ws.Range("E39").Value = "  -1.90%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01710"
$ws.Range("E40").Value = "  -1.74%  "
